$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.833.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.416.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "656.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.428"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.99%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.05"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.51%  "

# Row 10
$ws.Range("E10").Value = "  -0.04%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.408.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.33%  "

# Row 12
$ws.Range("E12").Value = "  +6.58%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +17.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000259"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.70%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.558.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.044.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +34.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.415.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.42%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +14.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.511"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +58.93%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +17.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "509.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000206"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.71%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.606.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.61%  "

# Row 30
$ws.Range("E30").Value = "  +13.03%  "

# Row 31
$ws.Range("E31").Value = "  +11.48%  "

# Row 32
$ws.Range("E32").Value = "  +4.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("B34").Value = "PolygonEcosystemToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.572"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +20.33%  "

# Row 35
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.54%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.64%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.156"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.89%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "514.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.33%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.858"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.51%  "

# Row 44
$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.30%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0419"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +23.64%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +16.40%  "

# Row 47
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.13%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.44%  "

# Row 49
$ws.Range("E49").Value = "  +0.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +16.62%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.60%  "
